$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number formats from column F (old column D) across to new D:E so the
# new cells inherit the same per-row style (date format row 7/38/80, numeric rows, etc.)
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D = 31-Dec-2018, E = 30-Sep-2018)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 386600
$ws.Range("E8").Value = 368700
$ws.Range("D9").Value = 332800
$ws.Range("E9").Value = 343400
$ws.Range("D10").Value = 53800
$ws.Range("E10").Value = 25300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 358800
$ws.Range("E17").Value = 361500
$ws.Range("D18").Value = 27800
$ws.Range("E18").Value = 7200
$ws.Range("D20").Value = -900
$ws.Range("E20").Value = -1500
$ws.Range("D21").Value = 41200
$ws.Range("E21").Value = 18700
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 26900
$ws.Range("E23").Value = 5700
$ws.Range("D24").Value = 6200
$ws.Range("E24").Value = 1800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 20700
$ws.Range("E26").Value = 3900
$ws.Range("D27").Value = 20700
$ws.Range("E27").Value = 3900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 100
$ws.Range("E29").Value = 1600
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 900
$ws.Range("E32").Value = 1500
$ws.Range("D33").Value = 20800
$ws.Range("E33").Value = 5500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 20800
$ws.Range("E35").Value = 5500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 9800
$ws.Range("E41").Value = 20200
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 160300
$ws.Range("E43").Value = 149100
$ws.Range("D44").Value = 137200
$ws.Range("E44").Value = 115000
$ws.Range("D45").Value = 3800
$ws.Range("E45").Value = 3900
$ws.Range("D46").Value = 311100
$ws.Range("E46").Value = 288200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 672200
$ws.Range("E48").Value = 639700
$ws.Range("D49").Value = 15000
$ws.Range("E49").Value = 15000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 36300
$ws.Range("E52").Value = 37300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1034600
$ws.Range("E54").Value = 980300
$ws.Range("D57").Value = 229500
$ws.Range("E57").Value = 209200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 9400
$ws.Range("D59").Value = 55200
$ws.Range("E59").Value = 28600
$ws.Range("D60").Value = 284700
$ws.Range("E60").Value = 247200
$ws.Range("D61").Value = 200000
$ws.Range("E61").Value = 190600
$ws.Range("D62").Value = 129600
$ws.Range("E62").Value = 133600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 614300
$ws.Range("E66").Value = 571400
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 187800
$ws.Range("E72").Value = 167100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 420300
$ws.Range("E76").Value = 408800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 20800
$ws.Range("E81").Value = 5500
$ws.Range("D83").Value = 14300
$ws.Range("E83").Value = 13000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 45700
$ws.Range("E89").Value = 50500
$ws.Range("D91").Value = -36600
$ws.Range("E91").Value = -19200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -37900
$ws.Range("E94").Value = -19600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -18200
$ws.Range("E100").Value = -27400
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -10400
$ws.Range("E102").Value = 3500
